$d = $word.ActiveDocument

# --------------------------------------------------------------------------
# This document originally has 3 paragraphs:
#   1. "CMP73010 Assignment 1 document (2017)"                     (title, unchanged)
#   2. "This is the document you will need to change. ..."        (to be replaced)
#   3. "Add instructions for your tutor to pull ... achieved."     (contains the
#      _GoBack bookmark; to be trimmed down to just the bookmark)
#
# Target layout (per the diff) has 4 paragraphs:
#   1. title                                                       (unchanged)
#   2. "Github" / " name: " / "jackmcloughlin"  (with proofErr spell wrappers)
#   3. "After I have edited ... / barrywilks / ) will receive ... choses to."
#   4. empty paragraph containing only the _GoBack bookmark
#
# We use Range.InsertXML (the WordOpenXML flat-OPC fragment form) to splice in
# exact OOXML, including the w:proofErr spell-check markers, which can't be
# produced through Find/Replace. NOTE: inserting 3+ <w:p> paragraphs in a
# single InsertXML call on this runtime corrupts the rest of the document, so
# the change is applied in two separate InsertXML calls, each introducing at
# most two paragraphs.
# --------------------------------------------------------------------------

function New-FlatOpcPackage($bodyXml) {
    return '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' + `
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' + `
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' + `
        '<pkg:xmlData>' + `
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
        '<w:body>' + $bodyXml + '</w:body>' + `
        '</w:document>' + `
        '</pkg:xmlData></pkg:part></pkg:package>'
}

# Step 1: turn paragraph 2 into the two "Github name" / "After I have edited" paragraphs.
$para2 = $d.Paragraphs(2)
$range1 = $d.Range($para2.Range.Start, $para2.Range.End)

$body1 = '<w:p w:rsidR="005D17C9" w:rsidRDefault="005D17C9">' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>Github</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t xml:space="preserve"> name: </w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>jackmcloughlin</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '</w:p>' + `
    '<w:p w:rsidR="005D17C9" w:rsidRDefault="005D17C9">' + `
    '<w:r><w:t>After I have edited and uploaded this copy to the repository. The owner of the master repository(</w:t></w:r>' + `
    '<w:proofErr w:type="spellStart"/>' + `
    '<w:r><w:t>barrywilks</w:t></w:r>' + `
    '<w:proofErr w:type="spellEnd"/>' + `
    '<w:r><w:t>) will receive a pull merge request. The requested file will be compared to the current file and both files will be merged if the user choses to.</w:t></w:r>' + `
    '</w:p>'

$range1.InsertXML((New-FlatOpcPackage $body1))

# Step 2: the old paragraph 3 ("Add instructions...") is now paragraph 4; shrink it
# down to just the _GoBack bookmark, removing all of its text runs.
$para4 = $d.Paragraphs(4)
$range2 = $d.Range($para4.Range.Start, $para4.Range.End)

$body2 = '<w:p w:rsidR="005D17C9" w:rsidRDefault="005D17C9">' + `
    '<w:bookmarkStart w:id="0" w:name="_GoBack"/>' + `
    '<w:bookmarkEnd w:id="0"/>' + `
    '</w:p>'

$range2.InsertXML((New-FlatOpcPackage $body2))
